$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Grow the table "Tabela1" from A1:B4 to A1:B10 (user added 6 more rows of data)
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B10"))

# 2) Write the new / changed cell values
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 3
$ws.Range("A5").Value = 12
$ws.Range("B5").Value = 4
$ws.Range("A6").Value = 12
$ws.Range("B6").Value = 5
$ws.Range("A7").Value = 12
$ws.Range("B7").Value = 6
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = 7
$ws.Range("A9").Value = 12
$ws.Range("B9").Value = 8
$ws.Range("A10").Value = 12
$ws.Range("B10").Value = 9

# 3) Extend the hidden _FilterDatabase defined name to match the new table range
$fdb = $wb.Names.Item(1)
$fdb.RefersTo = "='lista 1 kolumna'!`$A`$1:`$B`$10"

# 4) Re-scope the conditional formatting rules that lived on the old A2:A4 / A5:A1048576
#    ranges so they track the grown table (A2:A10 data rows, A11:A1048576 below it).
$cfA = $ws.Range("A1:A1048576").FormatConditions

$ruleDupBelow = $cfA.Item(1)
$ruleDupBelow.ModifyAppliesToRange($ws.Range("A11:A1048576"))

$ruleGreater = $cfA.Item(2)
$ruleGreater.ModifyAppliesToRange($ws.Range("A2:A10"))
$ruleGreater.Priority = 11

$ruleDup = $cfA.Item(3)
$ruleDup.ModifyAppliesToRange($ws.Range("A2:A10"))
$ruleDup.Priority = 12

# 5) Restore the active selection to C4, as recorded in the saved sheet view
$ws.Range("C4").Select() | Out-Null
